$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old row 2 data (A2:F2) which had a bug (duplicate sheet instantiation issue)
$ws.Range("A2:F2").ClearContents()

# Write corrected data: 3 rows (2-4), 5 columns (A-E) each
$data = @(
    @("Ares", "a1", "A/C", 18, $true),
    @("Ares", "a2", "A/C", 18, $true),
    @("Ares", "l2", "Lâmpada", 0, $false)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
